# Update countries & provincias Spain
# Applies the 31-May-2020 19:05 data refresh to the "Pais" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 31 de Mayo de 2020 a las 19:05"

# --- Country stat refresh (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) --------------

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1824047
$ws.Range("C4").Value = 7227
$ws.Range("D4").Value = 536293
$ws.Range("E4").Value = 1182034
$ws.Range("G4").Value = 163
$ws.Range("H4").Value = 105720

# Italia (row 9)
$ws.Range("B9").Value = 233019
$ws.Range("C9").Value = 355
$ws.Range("D9").Value = 157507
$ws.Range("E9").Value = 42097
$ws.Range("G9").Value = 75
$ws.Range("H9").Value = 33415

# India (row 10)
$ws.Range("B10").Value = 190536
$ws.Range("C10").Value = 8709
$ws.Range("D10").Value = 91621
$ws.Range("E10").Value = 93509
$ws.Range("G10").Value = 221
$ws.Range("H10").Value = 5406

# Alemania (row 12)
$ws.Range("B12").Value = 183420
$ws.Range("C12").Value = 126
$ws.Range("E12").Value = 9618

# Canada (row 17)
$ws.Range("B17").Value = 90516
$ws.Range("C17").Value = 326
$ws.Range("D17").Value = 48560
$ws.Range("E17").Value = 34864
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = 7092

# Irlanda (row 37)
$ws.Range("B37").Value = 24990
$ws.Range("C37").Value = 61
$ws.Range("E37").Value = 1249
$ws.Range("G37").Value = 1
$ws.Range("H37").Value = 1652

# Nepal (row 102)
$ws.Range("B102").Value = 1572
$ws.Range("C102").Value = 171
$ws.Range("E102").Value = 1344

# Santo Tome y Principe (row 139)
$ws.Range("B139").Value = 483
$ws.Range("C139").Value = 4
$ws.Range("E139").Value = 403

# Yemen (row 152)
$ws.Range("B152").Value = 323
$ws.Range("C152").Value = 13
$ws.Range("D152").Value = 14
$ws.Range("E152").Value = 229
$ws.Range("G152").Value = 3
$ws.Range("H152").Value = 80

# --- Reordered small-country rows -------------------------------------
# "Santa Lucia" now sorts ahead of "Belice"; swap the two rows' data
# (country name + Casos activos + Muertes) so row 200 carries
# Santa Lucia's figures and row 201 carries Belice's.
$ws.Range("A200").Value = "Santa Lucia"
$ws.Range("D200").Value = 18
$ws.Range("H200").Value = 0

$ws.Range("A201").Value = "Belice"
$ws.Range("D201").Value = 16
$ws.Range("H201").Value = 2

# "Papua Nueva Guinea" now sorts ahead of "Islas Virgenes Britanicas";
# swap the two rows' data the same way.
$ws.Range("A213").Value = "Papua Nueva Guinea"
$ws.Range("D213").Value = 8
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Islas Virgenes Britanicas"
$ws.Range("D214").Value = 7
$ws.Range("H214").Value = 1
